$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'320.64"
$ws.Range("E2").Value = "'3.78%"

# Row 3
$ws.Range("D3").Value = "'41.41"
$ws.Range("E3").Value = "'1.16%"

# Row 4
$ws.Range("D4").Value = "'5.248"
$ws.Range("E4").Value = "'2.54%"

# Row 6
$ws.Range("D6").Value = "'1.739"
$ws.Range("E6").Value = "'8.35%"

# Row 7
$ws.Range("D7").Value = "'0.9450"
$ws.Range("E7").Value = "'3.91%"

# Row 9
$ws.Range("D9").Value = "'0.1262"
$ws.Range("E9").Value = "'-1.87%"

# Row 10
$ws.Range("D10").Value = "'0.1862"
$ws.Range("E10").Value = "'3.27%"

# Row 11
$ws.Range("D11").Value = "'0.09200"
$ws.Range("E11").Value = "'1.32%"

# Row 12
$ws.Range("D12").Value = "'0.04169"
$ws.Range("E12").Value = "'-4.05%"

# Row 13
$ws.Range("D13").Value = "'0.1052"

# Row 14
$ws.Range("D14").Value = "'0.001287"
$ws.Range("E14").Value = "'3.20%"

# Row 15
$ws.Range("D15").Value = "'0.005839"
$ws.Range("E15").Value = "'1.94%"

# Row 17
$ws.Range("D17").Value = "'3.350"
$ws.Range("E17").Value = "'-0.07%"

# Row 18
$ws.Range("D18").Value = "'4.338"
$ws.Range("E18").Value = "'1.15%"

# Row 19
$ws.Range("E19").Value = "'1.22%"

# Row 20
$ws.Range("D20").Value = "'8.415"
$ws.Range("E20").Value = "'21.92%"

# Row 21
$ws.Range("E21").Value = "'-2.79%"

# Row 22
$ws.Range("D22").Value = "'0.2731"
$ws.Range("E22").Value = "'-0.24%"

# Row 23
$ws.Range("D23").Value = "'0.04030"
$ws.Range("E23").Value = "'-0.59%"

# Row 25
$ws.Range("D25").Value = "'0.004128"
$ws.Range("E25").Value = "'1.63%"

# Row 26
$ws.Range("E26").Value = "'0.02%"

# Row 38
$ws.Range("E38").Value = "'5.58%"

# Row 39
$ws.Range("D39").Value = "'0.05346"
$ws.Range("E39").Value = "'2.28%"

# Row 40
$ws.Range("D40").Value = "'0.007787"
$ws.Range("E40").Value = "'-0.80%"

# Row 41
$ws.Range("E41").Value = "'1.10%"

# Row 42
$ws.Range("D42").Value = "'0.007035"
$ws.Range("E42").Value = "'3.32%"

# Row 43
$ws.Range("E43").Value = "'6.95%"

# Row 44
$ws.Range("D44").Value = "'0.008302"
$ws.Range("E44").Value = "'11.84%"

# Row 45
$ws.Range("D45").Value = "'0.3465"
$ws.Range("E45").Value = "'3.68%"

# Row 46
$ws.Range("D46").Value = "'0.00006704"
$ws.Range("E46").Value = "'-2.47%"

# Row 47
$ws.Range("E47").Value = "'0.01%"

# Row 48
$ws.Range("D48").Value = "'0.1984"
$ws.Range("E48").Value = "'30.08%"

# Row 49
$ws.Range("D49").Value = "'0.004207"
$ws.Range("E49").Value = "'40.10%"

# Row 50
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.01%"

# Row 51
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.01%"
